$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shape = $s.Shapes.Item(3)
$shape.Table.ApplyStyle("{5A927BFA-8D41-47FC-AF4B-AC7AF810DE52}")
